# Auto-generated script to apply 2023-03-10 crime data update
# across the violent-crime-full-year workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 1193
$ws.Range('J3').Value = 1267
$ws.Range('H4').Value = 1687
$ws.Range('I4').Value = 1755
$ws.Range('J4').Value = 278
$ws.Range('I5').Value = 714
$ws.Range('H6').Value = 7918
$ws.Range('J6').Value = 1685
$ws.Range('I7').Value = 26197
$ws.Range('J7').Value = 4516

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 51

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 50
$ws.Range('J6').Value = 51
$ws.Range('J7').Value = 153

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 20
$ws.Range('J7').Value = 55

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J6').Value = 55
$ws.Range('J7').Value = 167

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J3').Value = 26
$ws.Range('J7').Value = 114

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J6').Value = 44
$ws.Range('J7').Value = 129
$ws.Range('J8').Value = 285
$ws.Range('J9').Value = 29
$ws.Range('J11').Value = 56
$ws.Range('J13').Value = 6
$ws.Range('J15').Value = 50
$ws.Range('J16').Value = 15
$ws.Range('J18').Value = 62
$ws.Range('H22').Value = 63
$ws.Range('J27').Value = 25
$ws.Range('J29').Value = 252
$ws.Range('J33').Value = 187
$ws.Range('J37').Value = 153
$ws.Range('J42').Value = 186
$ws.Range('J44').Value = 36
$ws.Range('J47').Value = 39
$ws.Range('J48').Value = 29
$ws.Range('J49').Value = 25
$ws.Range('J52').Value = 99
$ws.Range('J53').Value = 45
$ws.Range('J54').Value = 89
$ws.Range('J55').Value = 57
$ws.Range('H63').Value = 238
$ws.Range('I63').Value = 187
$ws.Range('J63').Value = 22
$ws.Range('J65').Value = 114
$ws.Range('J67').Value = 167
$ws.Range('J68').Value = 12
$ws.Range('J76').Value = 75
$ws.Range('J78').Value = 60
$ws.Range('J79').Value = 134
$ws.Range('J83').Value = 104
$ws.Range('J85').Value = 191
$ws.Range('J88').Value = 34
$ws.Range('J89').Value = 51
$ws.Range('J93').Value = 19
$ws.Range('J94').Value = 31
$ws.Range('J97').Value = 26
$ws.Range('J98').Value = 32
$ws.Range('J99').Value = 55
$ws.Range('I101').Value = 26197
$ws.Range('J101').Value = 4516

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J4').Value = 4
$ws.Range('J6').Value = 31
$ws.Range('J7').Value = 104

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J3').Value = 51
$ws.Range('J7').Value = 187

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('J3').Value = 9
$ws.Range('J7').Value = 25

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J4').Value = 5
$ws.Range('J6').Value = 46
$ws.Range('J7').Value = 89

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J3').Value = 94
$ws.Range('J7').Value = 252

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J2').Value = 13
$ws.Range('J4').Value = 4
$ws.Range('J6').Value = 9
$ws.Range('J7').Value = 36

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J2').Value = 4
$ws.Range('J6').Value = 16
$ws.Range('J7').Value = 29

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J6').Value = 43
$ws.Range('J7').Value = 75

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 50
$ws.Range('J3').Value = 69
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 191

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J6').Value = 14
$ws.Range('J7').Value = 44

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 31
$ws.Range('J7').Value = 186

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('J2').Value = 1
$ws.Range('J6').Value = 6

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J6').Value = 17
$ws.Range('J7').Value = 60

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J2').Value = 14
$ws.Range('J7').Value = 57

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J3').Value = 46
$ws.Range('J7').Value = 134

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J3').Value = 8
$ws.Range('J7').Value = 62

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('J2').Value = 2
$ws.Range('J7').Value = 19

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J6').Value = 40
$ws.Range('J7').Value = 99

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J6').Value = 18
$ws.Range('J7').Value = 31

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('J3').Value = 9
$ws.Range('J6').Value = 21
$ws.Range('J7').Value = 39

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J2').Value = 16
$ws.Range('J7').Value = 50

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J6').Value = 17
$ws.Range('J7').Value = 32

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 56

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('J3').Value = 9
$ws.Range('J7').Value = 29

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J6').Value = 19
$ws.Range('J7').Value = 26

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J6').Value = 16
$ws.Range('J7').Value = 34

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 96
$ws.Range('J7').Value = 285

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J6').Value = 12
$ws.Range('J7').Value = 25

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('J6').Value = 3
$ws.Range('J7').Value = 12

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J2').Value = 6
$ws.Range('J7').Value = 45

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('H6').Value = 18
$ws.Range('H7').Value = 63

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J6').Value = 39
$ws.Range('J7').Value = 129

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('J2').Value = 3
$ws.Range('J7').Value = 15
